# Rename the AHB-diff header columns from the generic "_old"/"_new" suffixes
# to the concrete format-version suffixes ("_FV2210" / "_FV2304"), wrap the
# data range in a real Excel Table (ListObject) and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the header row (row 1) text. These cells feed the shared-strings
#    table, so writing the new text here is equivalent to renaming the
#    "Segmentname_old" -> "Segmentname_FV2210" (and "_new" -> "_FV2304")
#    shared strings used by the diff.
$headers = @(
    "Segmentname_FV2210",
    "Segmentgruppe_FV2210",
    "Segment_FV2210",
    "Datenelement_FV2210",
    "Segment ID_FV2210",
    "Code_FV2210",
    "Qualifier_FV2210",
    "Beschreibung_FV2210",
    "Bedingungsausdruck_FV2210",
    "Bedingung_FV2210",
    "diff",
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)

for ($i = 0; $i -lt $headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# 2) Turn the A1:U61 data range into an Excel table (adds xl/tables/table1.xml
#    and the matching tableParts reference on the worksheet).
$dataRange = $ws.Range("A1:U61")
$table = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$table.Name = "Table1"
# No named table style (matches the plain tableStyleInfo emitted by the
# source export pipeline).
$table.TableStyle = ""

# 3) Freeze the header row (split below row 1, frozen at A2).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
